$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Fill in cell values
$ws.Range("A1").Value = "UserNamePassword"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Set the selection to C7
$ws.Range("C7").Select()
